# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# The two sheets mirror the same data, but one row (F9) received a
# slightly different updated value on each sheet, so they are handled
# with separate maps.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - row => new F value
$sheet1Updates = @{
    8  = 1653
    9  = 6068
    10 = 475
    11 = 343
    12 = 276
    15 = 128
    16 = 5401
    18 = 1262
    20 = 108
    21 = 218
    23 = 252
    24 = 95
    29 = 70
    32 = 40
    33 = 55
    35 = 58
    36 = 59
}

# Sheet "全部类型" (fourth sheet) - row => new F value
$sheet4Updates = @{
    8  = 1653
    9  = 6069
    10 = 475
    11 = 343
    12 = 276
    15 = 128
    16 = 5401
    18 = 1262
    20 = 108
    21 = 218
    23 = 252
    24 = 95
    29 = 70
    32 = 40
    33 = 55
    35 = 58
    36 = 59
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
